$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2632.6667
$ws.Range("J40").Value = 3450
$ws.Range("L40").Value = 3450
$ws.Range("N40").Value = -3800
$ws.Range("H132").Value = 71747.336
$ws.Range("I132").Value = 82323.84
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 246971.52
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -244441.52
$ws.Range("N132").Value = -14060
$ws.Range("H138").Value = 2622.7058
$ws.Range("J138").Value = 2801.3948
$ws.Range("L138").Value = 8404.1844
$ws.Range("N138").Value = -18684.1844
$ws.Range("H141").Value = 2372.6316
$ws.Range("I141").Value = 1828.2354
$ws.Range("K141").Value = 5484.706200000001
$ws.Range("M141").Value = -304.7062000000005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15940.109
$ws.Range("I32").Value = 16880.516
$ws.Range("K32").Value = 16880.516
$ws.Range("M32").Value = -16593.516
$ws.Range("H61").Value = 5569.5264
$ws.Range("I61").Value = 6858.2856
$ws.Range("J61").Value = 4817.75
$ws.Range("K61").Value = 6858.2856
$ws.Range("L61").Value = 4817.75
$ws.Range("M61").Value = -6646.2856
$ws.Range("N61").Value = -5241.75
$ws.Range("H63").Value = 3127001
$ws.Range("I63").Value = 2223.2222
$ws.Range("K63").Value = 2223.2222
$ws.Range("M63").Value = -1537.2222
$ws.Range("H66").Value = 3127001
$ws.Range("I66").Value = 2223.2222
$ws.Range("K66").Value = 11116.111
$ws.Range("M66").Value = -7684.111000000001
$ws.Range("H108").Value = 31925
$ws.Range("J108").Value = 31925
$ws.Range("L108").Value = 31925
$ws.Range("N108").Value = -39605
$ws.Range("H110").Value = 1953.2667
$ws.Range("I110").Value = 1811
$ws.Range("J110").Value = 2166.6667
$ws.Range("K110").Value = 1811
$ws.Range("L110").Value = 2166.6667
$ws.Range("M110").Value = 234
$ws.Range("N110").Value = -6256.6667
$ws.Range("H112").Value = 34191.145
$ws.Range("J112").Value = 34191.145
$ws.Range("L112").Value = 34191.145
$ws.Range("N112").Value = -37145.145
$ws.Range("H132").Value = 49356.727
$ws.Range("I132").Value = 4378
$ws.Range("J132").Value = 75058.86
$ws.Range("K132").Value = 13134
$ws.Range("L132").Value = 225176.58
$ws.Range("M132").Value = -10604
$ws.Range("N132").Value = -230236.58
$ws.Range("H136").Value = 5569.5264
$ws.Range("I136").Value = 6858.2856
$ws.Range("J136").Value = 4817.75
$ws.Range("K136").Value = 20574.8568
$ws.Range("L136").Value = 14453.25
$ws.Range("M136").Value = -18024.8568
$ws.Range("N136").Value = -19553.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 15745
$ws.Range("I82").Value = 5864.375
$ws.Range("K82").Value = 5864.375
$ws.Range("M82").Value = -5481.375
$ws.Range("H85").Value = 15745
$ws.Range("I85").Value = 5864.375
$ws.Range("K85").Value = 5864.375
$ws.Range("M85").Value = -4538.375
$ws.Range("H105").Value = 2001501.6
$ws.Range("I105").Value = 1402.5
$ws.Range("J105").Value = 3847747
$ws.Range("K105").Value = 1402.5
$ws.Range("L105").Value = 3847747
$ws.Range("M105").Value = 344.5
$ws.Range("N105").Value = -3851241
$ws.Range("H107").Value = 1570.9
$ws.Range("I107").Value = 1582.1111
$ws.Range("K107").Value = 1582.1111
$ws.Range("M107").Value = 337.8888999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 16348.967
$ws.Range("I31").Value = 37547.727
$ws.Range("J31").Value = 4076
$ws.Range("K31").Value = 37547.727
$ws.Range("L31").Value = 4076
$ws.Range("M31").Value = -37252.727
$ws.Range("N31").Value = -4666
$ws.Range("H34").Value = 16348.967
$ws.Range("I34").Value = 37547.727
$ws.Range("J34").Value = 4076
$ws.Range("K34").Value = 37547.727
$ws.Range("L34").Value = 4076
$ws.Range("M34").Value = -37345.727
$ws.Range("N34").Value = -4480
$ws.Range("H94").Value = 4494.154
$ws.Range("I94").Value = 2450
$ws.Range("K94").Value = 2450
$ws.Range("M94").Value = -1999
$ws.Range("H132").Value = 24451.791
$ws.Range("I132").Value = 28490.37
$ws.Range("J132").Value = 9105.200000000001
$ws.Range("K132").Value = 85471.11
$ws.Range("L132").Value = 27315.6
$ws.Range("M132").Value = -82941.11
$ws.Range("N132").Value = -32375.6
$ws.Range("H134").Value = 1009.2273
$ws.Range("I134").Value = 820.15
$ws.Range("J134").Value = 2900
$ws.Range("K134").Value = 2460.45
$ws.Range("L134").Value = 8700
$ws.Range("M134").Value = 74.55000000000018
$ws.Range("N134").Value = -13770

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 6729
$ws.Range("J68").Value = 9689.416999999999
$ws.Range("L68").Value = 29068.251
$ws.Range("N68").Value = -30690.251
$ws.Range("H71").Value = 6729
$ws.Range("J71").Value = 9689.416999999999
$ws.Range("L71").Value = 87204.753
$ws.Range("N71").Value = -95316.753
$ws.Range("H80").Value = 2933.3333
$ws.Range("J80").Value = 2933.3333
$ws.Range("L80").Value = 8799.999899999999
$ws.Range("N80").Value = -10671.9999
$ws.Range("H83").Value = 2933.3333
$ws.Range("J83").Value = 2933.3333
$ws.Range("L83").Value = 26399.9997
$ws.Range("N83").Value = -35759.9997
$ws.Range("H107").Value = 5445.6816
$ws.Range("J107").Value = 1047.4117
$ws.Range("L107").Value = 3142.2351
$ws.Range("N107").Value = -6982.2351
$ws.Range("H112").Value = 1783.3334
$ws.Range("I112").Value = 975
$ws.Range("K112").Value = 2925
$ws.Range("M112").Value = -1817
$ws.Range("H131").Value = 134152.94
$ws.Range("J131").Value = 150087.17
$ws.Range("L131").Value = 450261.51
$ws.Range("N131").Value = -460341.51
$ws.Range("H132").Value = 846.0833
$ws.Range("I132").Value = 728.1111
$ws.Range("J132").Value = 1200
$ws.Range("K132").Value = 6552.9999
$ws.Range("L132").Value = 10800
$ws.Range("M132").Value = -4022.9999
$ws.Range("N132").Value = -15860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 6325
$ws.Range("I5").Value = 500
$ws.Range("J5").Value = 8266.666999999999
$ws.Range("K5").Value = 500
$ws.Range("L5").Value = 8266.666999999999
$ws.Range("M5").Value = -388
$ws.Range("N5").Value = -8490.666999999999
$ws.Range("H102").Value = 1901.1111
$ws.Range("J102").Value = 2100.25
$ws.Range("L102").Value = 2100.25
$ws.Range("N102").Value = -5344.25
$ws.Range("H122").Value = 9500
$ws.Range("I122").Value = 9500
$ws.Range("K122").Value = 28500
$ws.Range("M122").Value = -26050
$ws.Range("H132").Value = 105813.266
$ws.Range("I132").Value = 106420.1
$ws.Range("J132").Value = 104599.6
$ws.Range("K132").Value = 319260.3
$ws.Range("L132").Value = 313798.8
$ws.Range("M132").Value = -316730.3
$ws.Range("N132").Value = -318858.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5053.478
$ws.Range("I7").Value = 4911.524
$ws.Range("J7").Value = 6544
$ws.Range("K7").Value = 4911.524
$ws.Range("L7").Value = 6544
$ws.Range("M7").Value = -4799.524
$ws.Range("N7").Value = -6768
$ws.Range("H40").Value = 114617
$ws.Range("I40").Value = 187648.33
$ws.Range("J40").Value = 5070
$ws.Range("K40").Value = 187648.33
$ws.Range("L40").Value = 5070
$ws.Range("M40").Value = -187512.33
$ws.Range("N40").Value = -5342
$ws.Range("H46").Value = 1098.4445
$ws.Range("I46").Value = 784.8
$ws.Range("K46").Value = 784.8
$ws.Range("M46").Value = -596.8
$ws.Range("H122").Value = 3406.2778
$ws.Range("I122").Value = 2810.3
$ws.Range("J122").Value = 4151.25
$ws.Range("K122").Value = 8430.900000000001
$ws.Range("L122").Value = 12453.75
$ws.Range("M122").Value = -5980.900000000001
$ws.Range("N122").Value = -17353.75
$ws.Range("H126").Value = 5053.478
$ws.Range("I126").Value = 4911.524
$ws.Range("J126").Value = 6544
$ws.Range("K126").Value = 14734.572
$ws.Range("L126").Value = 19632
$ws.Range("M126").Value = -12264.572
$ws.Range("N126").Value = -24572
$ws.Range("H132").Value = 3200.6667
$ws.Range("I132").Value = 2650
$ws.Range("J132").Value = 3310.8
$ws.Range("K132").Value = 7950
$ws.Range("L132").Value = 9932.400000000001
$ws.Range("M132").Value = -5420
$ws.Range("N132").Value = -14992.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1889.5714
$ws.Range("I107").Value = 541.4
$ws.Range("J107").Value = 2638.5557
$ws.Range("K107").Value = 1624.2
$ws.Range("L107").Value = 7915.6671
$ws.Range("M107").Value = 295.8000000000002
$ws.Range("N107").Value = -11755.6671
$ws.Range("H113").Value = 890.375
$ws.Range("I113").Value = 936.2
$ws.Range("J113").Value = 203
$ws.Range("K113").Value = 2808.6
$ws.Range("L113").Value = 609
$ws.Range("M113").Value = -638.6000000000004
$ws.Range("N113").Value = -4949
$ws.Range("H126").Value = 1870.7142
$ws.Range("I126").Value = 1559
$ws.Range("K126").Value = 4677
$ws.Range("M126").Value = -2207
$ws.Range("H132").Value = 2575.348
$ws.Range("I132").Value = 2770.8
$ws.Range("J132").Value = 2425
$ws.Range("K132").Value = 8312.400000000001
$ws.Range("L132").Value = 7275
$ws.Range("M132").Value = -5782.400000000001
$ws.Range("N132").Value = -12335
$ws.Range("H136").Value = 1796.0625
$ws.Range("I136").Value = 1231.3334
$ws.Range("J136").Value = 2134.9
$ws.Range("K136").Value = 3694.0002
$ws.Range("L136").Value = 6404.700000000001
$ws.Range("M136").Value = -1144.0002
$ws.Range("N136").Value = -11504.7
